$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of COVID data to append (dates 2022-02-07 .. 2022-02-13)
$rows = @(
    @("2022-02-07", "overview", "K02000001", "United Kingdom", 17866632, 57623, 45, 158363),
    @("2022-02-08", "overview", "K02000001", "United Kingdom", 17932803, 66183, 314, 158677),
    @("2022-02-09", "overview", "K02000001", "United Kingdom", 18000119, 68214, 276, 158953),
    @("2022-02-10", "overview", "K02000001", "United Kingdom", 18162199, 66638, 206, 159158),
    @("2022-02-11", "overview", "K02000001", "United Kingdom", 18220515, 58899, 193, 159351),
    @("2022-02-12", "overview", "K02000001", "United Kingdom", 18266015, 46025, 167, 159518),
    @("2022-02-13", "overview", "K02000001", "United Kingdom", 18306859, 41270, 52, 159570)
)

$startRow = 545

# Format column A for the new rows as Text first so the date-like strings
# ("YYYY-MM-DD") are kept as literal text instead of being auto-converted
# into Excel date serial numbers, matching the source data (inline strings).
$endRow = $startRow + $rows.Count - 1
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
